{"js": "// Add a new bullet list item \"Usage of internal keyword.\" after the last\n// paragraph in the document body (\"Physics2D Settings Unchecking Queries\n// Start in Colliders to disable detecting itself.\"), matching the same\n// ListParagraph / numbered-list formatting.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert a new paragraph right after the last one; Word carries over the\n// source paragraph's formatting (style + list numbering) automatically,\n// matching how a user pressing Enter at the end of a list item behaves.\nconst newParagraph = lastParagraph.insertParagraph(\"Usage of internal keyword.\", \"After\");\n\nawait context.sync();\n", "ps1": "# Add a new bullet list item \"Usage of internal keyword.\" after the last\n# paragraph in the document body (\"Physics2D Settings Unchecking Queries\n# Start in Colliders to disable detecting itself.\"), matching the same\n# ListParagraph / numbered-list formatting.\n\n$d = $word.ActiveDocument\n\n# Locate the last paragraph in the document (the final bullet item).\n$lastParagraph = $d.Paragraphs.Last\n\n# Insert a new paragraph right after it; Word carries over the source\n# paragraph's formatting (style + list numbering) automatically, matching\n# how a user pressing Enter at the end of a list item behaves.\n$lastParagraph.Range.InsertParagraphAfter()\n\n# Set the text of the newly created (now last) paragraph.\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"Usage of internal keyword.\"\n"}
